# Week 6 Meeting deck — "added timelog file to repo"
#
# PowerPoint had re-cached the auto-updating date fields (type="datetime1")
# that live on the slide master and every slide layout's Date Placeholder.
# Their displayed text moved from "11/12/2021" to "12/11/2021" (the field
# still auto-updates to "today" each time the deck is opened/saved — this
# just refreshes the cached text that ships in the file).
#
# Walk the slide master's own Date Placeholder plus every custom layout's
# Date Placeholder and re-stamp the cached text to the new value.

$p = $ppt.ActivePresentation
$oldDate = "11/12/2021"
$newDate = "12/11/2021"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
